$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pedidos")
$ws.Range("A11").Value = "a) pediu 5 potes de sorvetes"
$ws.Range("A12").Value = "b) foram 3, 2 de chocolate preto e 1 de chocolate branco"
